$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Prix Spot" - add column AY (03-aug) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Range("AX1").Copy() | Out-Null
$ws1.Range("AY1").PasteSpecial(-4122) | Out-Null
$ws1.Range("AY1").Value = "03-aug"

$ay1Values = @(
    80.40000000000001,
    69.26000000000001,
    42.96,
    37.38,
    30.02,
    30.2,
    32.03,
    26.28,
    6.16,
    0,
    -1.01,
    -2.1,
    -3,
    -9.9,
    -10.08,
    -4.98,
    -1.49,
    -0.01,
    12.51,
    40,
    78.2,
    96.13,
    97.40000000000001,
    85.09
)

for ($i = 0; $i -lt $ay1Values.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 51).Value = $ay1Values[$i]
}

# --- Sheet 2: "Gaz" - add row 48 ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A48").NumberFormat = "@"
$ws2.Range("A48").Value = "2025-08-01"
$ws2.Range("A48").Style = "Normal"
$ws2.Range("B48").Value = 32.65

# --- Sheet 3: "CO2" - add row 48 ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A48").NumberFormat = "@"
$ws3.Range("A48").Value = "2025-08-01"
$ws3.Range("A48").Style = "Normal"
$ws3.Range("B48").Value = 70.58
